# Add translation settings to the ODK-X "settings" sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("settings")

# --- Column widths (existing 3 cols get slightly wider, 5 new cols added) ---
# ColumnWidth input chosen so the underlying xlsx <col width> (which Excel
# rounds to whole pixels, width = (round(ColumnWidth*6)+5)/6) lands as close
# as possible to the target widths.
$ws.Columns.Item(1).ColumnWidth = 13.46
$ws.Columns.Item(2).ColumnWidth = 15.77
$ws.Columns.Item(3).ColumnWidth = 14.56
$ws.Columns.Item(4).ColumnWidth = 16.62
$ws.Columns.Item(5).ColumnWidth = 17.1
$ws.Columns.Item(6).ColumnWidth = 15.89
$ws.Columns.Item(7).ColumnWidth = 18.21
$ws.Columns.Item(8).ColumnWidth = 18.69

# --- Row 1: header row gets 5 new trailing columns (D..H) ---
$ws.Range("A1").Value = "setting_name"
$ws.Range("B1").Value = "value"
$ws.Range("C1").Value = "display.title.text"
$ws.Range("D1").Value = "display.title.text.pt"
$ws.Range("E1").Value = "display.title.text.sw"
$ws.Range("F1").Value = "display.locale.text"
$ws.Range("G1").Value = "display.locale.text.pt"
$ws.Range("H1").Value = "display.locale.text.sw"
$ws.Range("C1:H1").Style = "Normal"
$ws.Rows.Item(1).RowHeight = 13.8

# --- Row 2: form_id / hh_member_snake (unchanged content) ---
$ws.Range("A2").Value = "form_id"
$ws.Range("B2").Value = "hh_member_snake"

# --- Row 3: form_version (unchanged content) ---
$ws.Range("A3").Value = "form_version"

# --- Row 4: table_id / hh_member (unchanged content) ---
$ws.Range("A4").Value = "table_id"
$ws.Range("B4").Value = "hh_member"

# --- Row 5: survey / title text repeated across C,D,E; row height shrinks ---
$ws.Range("A5").Value = "survey"
$ws.Range("C5").Value = "Snake Bites Info"
$ws.Range("D5").Value = "Snake Bites Info"
$ws.Range("E5").Value = "Snake Bites Info"
$ws.Rows.Item(5).RowHeight = 12.8

# --- Row 6: instance_name / name (unchanged content) ---
$ws.Range("A6").Value = "instance_name"
$ws.Range("B6").Value = "name"

# --- Row 7: default locale = English ---
$ws.Range("A7").Value = "default"
$ws.Range("F7").Value = "English"
$ws.Range("G7").Value = "English"
$ws.Range("H7").Value = "English"
$ws.Range("A7:H7").Style = "Normal"
$ws.Rows.Item(7).RowHeight = 13.8

# --- Row 8: pt locale = Portuguese ---
$ws.Range("A8").Value = "pt"
$ws.Range("F8").Value = "Português"
$ws.Range("G8").Value = "Português"
$ws.Range("H8").Value = "Português"
$ws.Range("A8:H8").Style = "Normal"
$ws.Rows.Item(8).RowHeight = 13.8

# --- Row 9: sw locale = Swahili ---
$ws.Range("A9").Value = "sw"
$ws.Range("F9").Value = "Kiswahili"
$ws.Range("G9").Value = "Kiswahili"
$ws.Range("H9").Value = "Kiswahili"
$ws.Range("A9:H9").Style = "Normal"
$ws.Rows.Item(9).RowHeight = 13.8
